{"js": "// Apply the text replacements described by the diff: the date line and\n// the 25 \"NNN\u00d7N=\" multiplication prompts scattered through the table.\n// Every \"before\" string is unique in the document, so a simple\n// search-and-replace per pair is safe and keeps existing run formatting\n// (fonts/size) intact because we replace only the found sub-range.\nconst replacements = [\n  [\"2025-02-01 Saturday\", \"2025-02-02 Sunday\"],\n  [\"793\u00d75=\", \"491\u00d74=\"],\n  [\"849\u00d77=\", \"192\u00d76=\"],\n  [\"285\u00d74=\", \"770\u00d78=\"],\n  [\"817\u00d72=\", \"184\u00d76=\"],\n  [\"816\u00d73=\", \"945\u00d78=\"],\n  [\"905\u00d76=\", \"296\u00d72=\"],\n  [\"642\u00d73=\", \"906\u00d74=\"],\n  [\"619\u00d79=\", \"497\u00d74=\"],\n  [\"255\u00d76=\", \"243\u00d74=\"],\n  [\"550\u00d73=\", \"742\u00d76=\"],\n  [\"719\u00d76=\", \"111\u00d72=\"],\n  [\"108\u00d75=\", \"439\u00d76=\"],\n  [\"749\u00d77=\", \"961\u00d78=\"],\n  [\"527\u00d73=\", \"461\u00d72=\"],\n  [\"476\u00d78=\", \"329\u00d79=\"],\n  [\"626\u00d79=\", \"691\u00d76=\"],\n  [\"877\u00d78=\", \"204\u00d73=\"],\n  [\"395\u00d74=\", \"848\u00d77=\"],\n  [\"809\u00d75=\", \"332\u00d74=\"],\n  [\"690\u00d72=\", \"333\u00d78=\"],\n  [\"672\u00d74=\", \"216\u00d72=\"],\n  [\"888\u00d78=\", \"419\u00d75=\"],\n  [\"797\u00d77=\", \"417\u00d77=\"],\n  [\"614\u00d79=\", \"991\u00d74=\"],\n  [\"555\u00d75=\", \"905\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: the date line and\n# the 25 \"NNN\u00d7N=\" multiplication prompts scattered through the table.\n# Every \"before\" string is unique in the document, so a simple\n# Find/Replace (wdReplaceAll) per pair is safe and keeps existing run\n# formatting (fonts/size) intact since Find/Replace only touches the\n# matched text, not the surrounding run properties.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old = \"2025-02-01 Saturday\"; new = \"2025-02-02 Sunday\"},\n  @{old = \"793\u00d75=\"; new = \"491\u00d74=\"},\n  @{old = \"849\u00d77=\"; new = \"192\u00d76=\"},\n  @{old = \"285\u00d74=\"; new = \"770\u00d78=\"},\n  @{old = \"817\u00d72=\"; new = \"184\u00d76=\"},\n  @{old = \"816\u00d73=\"; new = \"945\u00d78=\"},\n  @{old = \"905\u00d76=\"; new = \"296\u00d72=\"},\n  @{old = \"642\u00d73=\"; new = \"906\u00d74=\"},\n  @{old = \"619\u00d79=\"; new = \"497\u00d74=\"},\n  @{old = \"255\u00d76=\"; new = \"243\u00d74=\"},\n  @{old = \"550\u00d73=\"; new = \"742\u00d76=\"},\n  @{old = \"719\u00d76=\"; new = \"111\u00d72=\"},\n  @{old = \"108\u00d75=\"; new = \"439\u00d76=\"},\n  @{old = \"749\u00d77=\"; new = \"961\u00d78=\"},\n  @{old = \"527\u00d73=\"; new = \"461\u00d72=\"},\n  @{old = \"476\u00d78=\"; new = \"329\u00d79=\"},\n  @{old = \"626\u00d79=\"; new = \"691\u00d76=\"},\n  @{old = \"877\u00d78=\"; new = \"204\u00d73=\"},\n  @{old = \"395\u00d74=\"; new = \"848\u00d77=\"},\n  @{old = \"809\u00d75=\"; new = \"332\u00d74=\"},\n  @{old = \"690\u00d72=\"; new = \"333\u00d78=\"},\n  @{old = \"672\u00d74=\"; new = \"216\u00d72=\"},\n  @{old = \"888\u00d78=\"; new = \"419\u00d75=\"},\n  @{old = \"797\u00d77=\"; new = \"417\u00d77=\"},\n  @{old = \"614\u00d79=\"; new = \"991\u00d74=\"},\n  @{old = \"555\u00d75=\"; new = \"905\u00d79=\"}\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $p.old\n  $find.Replacement.Text = $p.new\n  $find.Execute($p.old, $false, $true, $false, $false, $false, $true, 1, $false, $p.new, 2) | Out-Null\n}\n"}
